$wb = $excel.ActiveWorkbook

# ===== ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4247.5654
$ws.Range("I40").Value = 1934
$ws.Range("K40").Value = 1934
$ws.Range("M40").Value = -1759
$ws.Range("H64").Value = 7058.16
$ws.Range("I64").Value = 5144.5454
$ws.Range("J64").Value = 8561.714
$ws.Range("K64").Value = 5144.5454
$ws.Range("L64").Value = 8561.714
$ws.Range("M64").Value = -4896.5454
$ws.Range("N64").Value = -9057.714
$ws.Range("H67").Value = 7058.16
$ws.Range("I67").Value = 5144.5454
$ws.Range("J67").Value = 8561.714
$ws.Range("K67").Value = 5144.5454
$ws.Range("L67").Value = 8561.714
$ws.Range("M67").Value = -4286.5454
$ws.Range("N67").Value = -10277.714
$ws.Range("H98").Value = 1095.5
$ws.Range("I98").Value = 1150.05
$ws.Range("J98").Value = 550
$ws.Range("K98").Value = 1150.05
$ws.Range("L98").Value = 550
$ws.Range("M98").Value = 347.95
$ws.Range("N98").Value = -3546
$ws.Range("H107").Value = 66668170
$ws.Range("I107").Value = 66668170
$ws.Range("K107").Value = 66668170
$ws.Range("M107").Value = -66666250
$ws.Range("H122").Value = 1095.5
$ws.Range("I122").Value = 1150.05
$ws.Range("J122").Value = 550
$ws.Range("K122").Value = 3450.15
$ws.Range("L122").Value = 1650
$ws.Range("M122").Value = -1000.15
$ws.Range("N122").Value = -6550
$ws.Range("H132").Value = 30306366
$ws.Range("I132").Value = 47623090
$ws.Range("K132").Value = 142869270
$ws.Range("M132").Value = -142866740
$ws.Range("H135").Value = 838.275
$ws.Range("I135").Value = 810.9697
$ws.Range("K135").Value = 7298.7273
$ws.Range("M135").Value = -4763.7273
$ws.Range("H137").Value = 68669.89
$ws.Range("I137").Value = 95688.78999999999
$ws.Range("K137").Value = 287066.37
$ws.Range("M137").Value = -284516.37
$ws.Range("H141").Value = 1554.4
$ws.Range("I141").Value = 1605.5
$ws.Range("J141").Value = 1350
$ws.Range("K141").Value = 4816.5
$ws.Range("L141").Value = 4050
$ws.Range("M141").Value = 363.5
$ws.Range("N141").Value = -14410

# ===== ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3286.388
$ws.Range("I32").Value = 1929.8334
$ws.Range("J32").Value = 8921.308000000001
$ws.Range("K32").Value = 1929.8334
$ws.Range("L32").Value = 8921.308000000001
$ws.Range("M32").Value = -1642.8334
$ws.Range("N32").Value = -9495.308000000001
$ws.Range("H61").Value = 3377.7856
$ws.Range("I61").Value = 2867.4443
$ws.Range("K61").Value = 2867.4443
$ws.Range("M61").Value = -2655.4443
$ws.Range("H74").Value = 53128.42
$ws.Range("I74").Value = 4000.1724
$ws.Range("J74").Value = 211430.56
$ws.Range("K74").Value = 4000.1724
$ws.Range("L74").Value = 211430.56
$ws.Range("M74").Value = -3126.1724
$ws.Range("N74").Value = -213178.56
$ws.Range("H77").Value = 53128.42
$ws.Range("I77").Value = 4000.1724
$ws.Range("J77").Value = 211430.56
$ws.Range("K77").Value = 20000.862
$ws.Range("L77").Value = 1057152.8
$ws.Range("M77").Value = -15632.862
$ws.Range("N77").Value = -1065888.8
$ws.Range("H82").Value = 62582
$ws.Range("J82").Value = 65000
$ws.Range("L82").Value = 65000
$ws.Range("N82").Value = -65722
$ws.Range("H85").Value = 62582
$ws.Range("J85").Value = 65000
$ws.Range("L85").Value = 65000
$ws.Range("N85").Value = -67496
$ws.Range("H132").Value = 3168.92
$ws.Range("I132").Value = 2181.6875
$ws.Range("J132").Value = 4924
$ws.Range("K132").Value = 6545.0625
$ws.Range("L132").Value = 14772
$ws.Range("M132").Value = -4015.0625
$ws.Range("N132").Value = -19832
$ws.Range("H136").Value = 3377.7856
$ws.Range("I136").Value = 2867.4443
$ws.Range("K136").Value = 8602.332900000001
$ws.Range("M136").Value = -6052.332900000001

# ===== BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 3214
$ws.Range("I75").Value = 3214
$ws.Range("K75").Value = 3214
$ws.Range("M75").Value = -2278
$ws.Range("H78").Value = 3214
$ws.Range("I78").Value = 3214
$ws.Range("K78").Value = 9642
$ws.Range("M78").Value = -4962
$ws.Range("H82").Value = 16239.777
$ws.Range("I82").Value = 9519.875
$ws.Range("K82").Value = 9519.875
$ws.Range("M82").Value = -9136.875
$ws.Range("H85").Value = 16239.777
$ws.Range("I85").Value = 9519.875
$ws.Range("K85").Value = 9519.875
$ws.Range("M85").Value = -8193.875
$ws.Range("H86").Value = 3708578.2
$ws.Range("I86").Value = 4767591
$ws.Range("J86").Value = 2034
$ws.Range("K86").Value = 4767591
$ws.Range("L86").Value = 2034
$ws.Range("M86").Value = -4766468
$ws.Range("N86").Value = -4280
$ws.Range("H89").Value = 3708578.2
$ws.Range("I89").Value = 4767591
$ws.Range("J89").Value = 2034
$ws.Range("K89").Value = 23837955
$ws.Range("L89").Value = 10170
$ws.Range("M89").Value = -23832339
$ws.Range("N89").Value = -21402
$ws.Range("H94").Value = 3966007.2
$ws.Range("I94").Value = 6994192
$ws.Range("K94").Value = 6994192
$ws.Range("M94").Value = -6993741

# ===== CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H47").Value = 20000
$ws.Range("I47").Value = 20000
$ws.Range("K47").Value = 20000
$ws.Range("M47").Value = -19434
$ws.Range("H58").Value = 2482.4194
$ws.Range("I58").Value = 2042.6666
$ws.Range("J58").Value = 3405.9
$ws.Range("K58").Value = 2042.6666
$ws.Range("L58").Value = 3405.9
$ws.Range("M58").Value = -1839.6666
$ws.Range("N58").Value = -3811.9
$ws.Range("H132").Value = 56109.17
$ws.Range("I132").Value = 40497.58
$ws.Range("K132").Value = 121492.74
$ws.Range("M132").Value = -118962.74
$ws.Range("H134").Value = 2978.3044
$ws.Range("I134").Value = 1357.6428
$ws.Range("J134").Value = 5499.3335
$ws.Range("K134").Value = 4072.9284
$ws.Range("L134").Value = 16498.0005
$ws.Range("M134").Value = -1537.9284
$ws.Range("N134").Value = -21568.0005
$ws.Range("H136").Value = 2482.4194
$ws.Range("I136").Value = 2042.6666
$ws.Range("J136").Value = 3405.9
$ws.Range("K136").Value = 6127.9998
$ws.Range("L136").Value = 10217.7
$ws.Range("M136").Value = -3577.9998
$ws.Range("N136").Value = -15317.7

# ===== CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 386.34482
$ws.Range("I2").Value = 194.375
$ws.Range("J2").Value = 459.4762
$ws.Range("K2").Value = 1166.25
$ws.Range("L2").Value = 2756.8572
$ws.Range("M2").Value = -1053.25
$ws.Range("N2").Value = -2982.8572
$ws.Range("H7").Value = 2540.125
$ws.Range("I7").Value = 2936.8333
$ws.Range("K7").Value = 8810.499899999999
$ws.Range("M7").Value = -8698.499899999999
$ws.Range("H11").Value = 12132.143
$ws.Range("I11").Value = 10012.5
$ws.Range("K11").Value = 30037.5
$ws.Range("M11").Value = -29897.5
$ws.Range("H35").Value = 425
$ws.Range("I35").Value = 433.33334
$ws.Range("J35").Value = 400
$ws.Range("K35").Value = 1300.00002
$ws.Range("L35").Value = 1200
$ws.Range("M35").Value = -1012.00002
$ws.Range("N35").Value = -1776
$ws.Range("H82").Value = 2806.5
$ws.Range("I82").Value = 2806.5
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 8419.5
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -8013.5
$ws.Range("H85").Value = 2806.5
$ws.Range("I85").Value = 2806.5
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 8419.5
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -7015.5
$ws.Range("H103").Value = 166.66667
$ws.Range("I103").Value = 166.66667
$ws.Range("K103").Value = 500.00001
$ws.Range("M103").Value = 378.99999
$ws.Range("H132").Value = 1634.7333
$ws.Range("J132").Value = 2338.6667
$ws.Range("L132").Value = 21048.0003
$ws.Range("N132").Value = -26108.0003

# ===== GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 983.8461
$ws.Range("I7").Value = 482.5
$ws.Range("K7").Value = 482.5
$ws.Range("M7").Value = -370.5
$ws.Range("H8").Value = 983.8461
$ws.Range("I8").Value = 482.5
$ws.Range("K8").Value = 482.5
$ws.Range("M8").Value = -343.5
$ws.Range("H9").Value = 2199.3333
$ws.Range("I9").Value = 732
$ws.Range("J9").Value = 3666.6667
$ws.Range("K9").Value = 732
$ws.Range("L9").Value = 3666.6667
$ws.Range("M9").Value = -562
$ws.Range("N9").Value = -4006.6667
$ws.Range("H70").Value = 8338422.5
$ws.Range("I70").Value = 10531126
$ws.Range("K70").Value = 10531126
$ws.Range("M70").Value = -10530856
$ws.Range("H73").Value = 8338422.5
$ws.Range("I73").Value = 10531126
$ws.Range("K73").Value = 10531126
$ws.Range("M73").Value = -10530190
$ws.Range("H80").Value = 1879385.8
$ws.Range("I80").Value = 2711767.5
$ws.Range("K80").Value = 2711767.5
$ws.Range("M80").Value = -2710769.5
$ws.Range("H83").Value = 1879385.8
$ws.Range("I83").Value = 2711767.5
$ws.Range("K83").Value = 13558837.5
$ws.Range("M83").Value = -13553845.5

# ===== LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 20004334
$ws.Range("I2").Value = 20004334
$ws.Range("K2").Value = 20004334
$ws.Range("M2").Value = -20004222
$ws.Range("H93").Value = 37039990
$ws.Range("I93").Value = 55557820
$ws.Range("K93").Value = 55557820
$ws.Range("M93").Value = -55556572
$ws.Range("H132").Value = 9985.923000000001
$ws.Range("I132").Value = 10382.5
$ws.Range("K132").Value = 31147.5
$ws.Range("M132").Value = -28617.5

# ===== WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 45000
$ws.Range("J110").Value = 45000
$ws.Range("L110").Value = 45000
$ws.Range("N110").Value = -53180
$ws.Range("H126").Value = 2793.5625
$ws.Range("I126").Value = 2481.6365
$ws.Range("J126").Value = 3479.8
$ws.Range("K126").Value = 7444.9095
$ws.Range("L126").Value = 10439.4
$ws.Range("M126").Value = -4974.9095
$ws.Range("N126").Value = -15379.4
$ws.Range("H136").Value = 5000.909
$ws.Range("I136").Value = 2960
$ws.Range("K136").Value = 8880
$ws.Range("M136").Value = -6330
